$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data (player, position, team) for rows 2-17 after the edit:
#  - "Kelly Oubre Jr." / "SG,SF" / "Philadelphia 76ers" is a new row inserted
#    right after "Dyson Daniels" (pushing the rows below it down by one).
#  - "Victor Wembanyama" moves up to sit right before "Alperen Sengün".
#  - "Deandre Ayton" / "Portland Trail Blazers" (the last row) is removed.
$data = @(
    @("Donte DiVincenzo",   "PG,SG,SF", "Minnesota Timberwolves"),
    @("Josh Hart",          "SG,SF,PF", "New York Knicks"),
    @("Dyson Daniels",      "PG,SG,SF", "Atlanta Hawks"),
    @("Kelly Oubre Jr.",    "SG,SF",    "Philadelphia 76ers"),
    @("Andrew Wiggins",     "SF,PF",    "Golden State Warriors"),
    @("Chris Boucher",      "PF,C",     "Toronto Raptors"),
    @("Michael Porter Jr.", "SF,PF",    "Denver Nuggets"),
    @("De'Andre Hunter",    "SF,PF",    "Atlanta Hawks"),
    @("Victor Wembanyama",  "C",        "San Antonio Spurs"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Santi Aldama",       "PF,C",     "Memphis Grizzlies"),
    @("Kristaps Porzingis", "PF,C",     "Boston Celtics"),
    @("Donovan Mitchell",   "PG,SG",    "Cleveland Cavaliers"),
    @("Domantas Sabonis",   "C",        "Sacramento Kings"),
    @("Malik Beasley",      "SG,SF",    "Detroit Pistons"),
    @("Cam Thomas",         "SG,SF",    "Brooklyn Nets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
